$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 21; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $cell.Value2 = $cell.Value2 + 154
    $ws.Rows.Item($r).RowHeight = 13.8
}

$ws.Range("C2:C21").Select()
